$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.014545559883118
$ws.Range("B1").Value = 3.117747783660889
$ws.Range("C1").Value = 6.682947158813477
$ws.Range("D1").Value = 1.890783429145813
$ws.Range("E1").Value = 1.329131603240967
